$d = $word.ActiveDocument

# 1) Merge "Name: " + "Jason Liu" into a single run's text "Name: Jason Liu"
$d.Content.Find.Execute("Name: Jason Liu", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Name: Jason Liu", 2)

# 2) Merge "Period: " + "5" into a single run's text "Period: 5"
$d.Content.Find.Execute("Period: 5", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Period: 5", 2)

# 3) Merge the two sentences about the best model / accuracy into one run
$old3 = "The best model we used was a combination of XOR and XNOR with an AND gate for creating a square model. Through the sigmoid function, this rounds out the edge to create a circle. The accuracy of this model was ___. "
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false,
                         $true, 1, $false, $old3, 2)
